# Automatische test-sync: 2025-08-03 14:36:50
# Append a new logged mail entry (row 13) to the "Logs" sheet and refresh
# the derived "Overig" count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 13

$ws.Cells.Item($newRow, 1).Value  = "Wil je dit oppakken?"
$ws.Cells.Item($newRow, 2).Value  = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 3).Value  = "Testmail #2: Wil je dit oppakken?"
$ws.Cells.Item($newRow, 4).Value  = "Overig"
$ws.Cells.Item($newRow, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$ws.Cells.Item($newRow, 6).Value  = "2025-08-03 14:35:56"
$ws.Cells.Item($newRow, 7).Value  = "Ja"
$ws.Cells.Item($newRow, 8).Value  = "Ja"
$ws.Cells.Item($newRow, 9).Value  = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# The sheet's conditional formatting rules were scoped to rows 2-12; grow
# each of them by one row so the newly appended row is covered too.
$ws.Range("D2:D12").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D13"))
$ws.Range("G2:G12").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G13"))
$ws.Range("H2:H12").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H13"))
$ws.Range("I2:I12").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I13"))
$ws.Range("J2:J12").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J13"))

# Update the Dashboard summary count for the "Overig" category to include
# the new row.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(4, 2).Value = 3
